# daily auto push: 2026-02-16 05:09 UTC
# Insert a new data row at row 830 (2026/02/16, 月, 13, 201), pushing the
# existing rows 830-871 down to 831-872, and extending the sheet's used
# range from D871 to D872.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 830..871 down by one, opening up a blank row 830.
$ws.Rows.Item(830).Insert()

# Fill the newly opened row. The date column stores plain text like
# "2026/12/29" elsewhere in the sheet (not a real date value), so force
# text with a leading apostrophe and then drop the resulting "Text" number
# format so the new cell's style matches its unstyled neighbours.
$ws.Range("A830").Value = "'2026/02/16"
$ws.Range("A830").ClearFormats()
$ws.Range("B830").Value = "月"
$ws.Range("C830").Value = 13
$ws.Range("D830").Value = 201
